# Updated code to use relative paths + added Excel summary files
# Refresh KPI_Summary values (rows 2-8) with the newly-computed statistics,
# and restore the correct KPI label ordering for rows 4-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - LPE
$ws.Range("A2").Value = 'LPE'
$ws.Range("B2").Value = [double]"0.8919753086419752"
$ws.Range("C2").Value = [double]"0.7212962962962963"
$ws.Range("D2").Value = [double]"0.6180555555555556"
$ws.Range("E2").Value = [double]"0.7100694444444444"
$ws.Range("F2").Value = [double]"-0.2739197530864197"
$ws.Range("G2").Value = [double]"-0.01122685185185185"
$ws.Range("H2").Value = [double]"-0.2626929012345679"
$ws.Range("I2").Value = [double]"-0.2626929012345677"
$ws.Range("J2").Value = [double]"-1.384608271639888"
$ws.Range("K2").Value = [double]"0.1984710570047823"

# Row 3 - avg_loan_size
$ws.Range("A3").Value = 'avg_loan_size'
$ws.Range("B3").Value = [double]"6981.406522119341"
$ws.Range("C3").Value = [double]"6853.072556040564"
$ws.Range("D3").Value = [double]"7069.619942332683"
$ws.Range("E3").Value = [double]"6777.255993234624"
$ws.Range("F3").Value = [double]"88.21342021334193"
$ws.Range("G3").Value = [double]"-235.3117629562824"
$ws.Range("H3").Value = [double]"323.5251831696243"
$ws.Range("I3").Value = [double]"164.0299830192826"
$ws.Range("J3").Value = [double]"0.2301040328135855"
$ws.Range("K3").Value = [double]"0.82290584471427"

# Row 4 - dq30_pct_unit
$ws.Range("A4").Value = 'dq30_pct_unit'
$ws.Range("B4").Value = [double]"0.03893627282589385"
$ws.Range("C4").Value = [double]"0.03601152563666175"
$ws.Range("D4").Value = [double]"0.0322940045920654"
$ws.Range("E4").Value = [double]"0.03024964109376913"
$ws.Range("F4").Value = [double]"-0.006642268233828448"
$ws.Range("G4").Value = [double]"-0.005761884542892614"
$ws.Range("H4").Value = [double]"-0.0008803836909358339"
$ws.Range("I4").Value = [double]"-0.0008803836909358209"
$ws.Range("J4").Value = [double]"-0.5544772572166765"
$ws.Range("K4").Value = [double]"0.5916177597107553"

# Row 5 - dq30_pct_$
$ws.Range("A5").Value = 'dq30_pct_$'
$ws.Range("B5").Value = [double]"0.9965966054932998"
$ws.Range("C5").Value = [double]"0.9982317537512463"
$ws.Range("D5").Value = [double]"0.9959517694192375"
$ws.Range("E5").Value = [double]"0.9976286479793253"
$ws.Range("F5").Value = [double]"-0.0006448360740621087"
$ws.Range("G5").Value = [double]"-0.0006031057719208243"
$ws.Range("H5").Value = [double]"-4.173030214128445e-05"
$ws.Range("I5").Value = [double]"-4.173030214127404e-05"
$ws.Range("J5").Value = [double]"-0.07934386342083714"
$ws.Range("K5").Value = [double]"0.9383720275384938"

# Row 6 - dq29_pot30_payment_rate_unit_per_day
$ws.Range("A6").Value = 'dq29_pot30_payment_rate_unit_per_day'
$ws.Range("B6").Value = [double]"0.007530871498455879"
$ws.Range("C6").Value = [double]"0.008874039540627269"
$ws.Range("D6").Value = [double]"0.009285787629030989"
$ws.Range("E6").Value = [double]"0.006929599208516626"
$ws.Range("F6").Value = [double]"0.001754916130575109"
$ws.Range("G6").Value = [double]"-0.001944440332110645"
$ws.Range("H6").Value = [double]"0.003699356462685754"
$ws.Range("I6").Value = [double]"0.003699356462685753"
$ws.Range("J6").Value = [double]"1.700028017452019"
$ws.Range("K6").Value = [double]"0.1209561116582869"

# Row 7 - dq29_pot30_payment_rate_unit_up_to_day
$ws.Range("A7").Value = 'dq29_pot30_payment_rate_unit_up_to_day'
$ws.Range("B7").Value = [double]"0.6214766727126475"
$ws.Range("C7").Value = [double]"0.579963746436155"
$ws.Range("D7").Value = [double]"0.4432950939350613"
$ws.Range("E7").Value = [double]"0.4068627846668502"
$ws.Range("F7").Value = [double]"-0.1781815787775861"
$ws.Range("G7").Value = [double]"-0.1731009617693049"
$ws.Range("H7").Value = [double]"-0.005080617008281246"
$ws.Range("I7").Value = [double]"-0.005080617008281385"
$ws.Range("J7").Value = [double]"-0.4656681151119703"
$ws.Range("K7").Value = [double]"0.6490915694341669"

# Row 8 - dq29_pot30_payment_rate_$_up_to_day
$ws.Range("A8").Value = 'dq29_pot30_payment_rate_$_up_to_day'
$ws.Range("B8").Value = [double]"0.002440572273129786"
$ws.Range("C8").Value = [double]"0.001103606986237561"
$ws.Range("D8").Value = [double]"0.004638147506546255"
$ws.Range("E8").Value = [double]"0.002730647807165762"
$ws.Range("F8").Value = [double]"0.002197575233416469"
$ws.Range("G8").Value = [double]"0.001627040820928201"
$ws.Range("H8").Value = [double]"0.0005705344124882672"
$ws.Range("I8").Value = [double]"0.0005705344124882679"
$ws.Range("J8").Value = [double]"0.7654840455173647"
$ws.Range("K8").Value = [double]"0.4562127228975335"

